$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "dsmljcnsd"
$ws.Range("A2").Value = "dcdsjk "
$ws.Range("A3").Value = "dclkjsdkcunjk n"
$ws.Range("A4").Value = "sd;clksdjcjn"
$ws.Range("A5").Value = "clk sjdk n"
$ws.Range("A6").Value = "sd snjij"

$ws.Range("A7").Select()
